$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.246.14"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.248.21"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.16"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.59"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.520"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.17"
$ws.Range("E10").Value = "  +5.70%  "
$ws.Range("E11").Value = "  +5.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0793"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.57"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.599.03"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.09"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.261.16"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.750"
$ws.Range("E18").Value = "  +1.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.150.82"
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.77"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.84"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.69"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.46"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.91"
$ws.Range("E28").Value = "  +4.24%  "
$ws.Range("E29").Value = "  -8.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.53"
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.45"
$ws.Range("E31").Value = "  +3.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.19"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.11"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0729"
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  +6.79%  "
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.33"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("E41").Value = "  +5.10%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.077.14"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.19"
$ws.Range("E44").Value = "  +10.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.39"
$ws.Range("E45").Value = "  +6.34%  "
$ws.Range("E46").Value = "  +2.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("E47").Value = "  +6.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.91"
$ws.Range("E48").Value = "  -10.00%  "
$ws.Range("E49").Value = "  +3.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.471.88"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("E51").Value = "  +4.48%  "
